# main_conditions_26.xlsx stimulus update
# - rename the "face" image category to "book" everywhere it is used as a
#   stimulus filename (promptFile / correctFile / dist_01File / dist_02File
#   columns A-D), e.g. "face//face_12.jpg" -> "book//book_12.jpg"
# - expand the abbreviated answer codes in column L (correct_ans) to their
#   full words: y -> left, r -> right, b -> center

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow   = $usedRange.Rows.Count
$lastCol   = $usedRange.Columns.Count

for ($row = 2; $row -le $lastRow; $row++) {

    # Columns A-D hold the stimulus file paths (promptFile, correctFile,
    # dist_01File, dist_02File). Swap the "face" category for "book" in any
    # of them, keeping the numeric suffix untouched.
    for ($col = 1; $col -le 4; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $value = $cell.Value2
        if ($value -ne $null -and $value -like "face//face_*") {
            $cell.Value2 = ($value -replace "face", "book")
        }
    }

    # Column L (correct_ans) uses single-letter shorthand for the answer
    # position - spell it out in full.
    $answerCell = $ws.Cells.Item($row, 12)
    $answer = $answerCell.Value2
    if ($answer -eq "y") {
        $answerCell.Value2 = "left"
    } elseif ($answer -eq "r") {
        $answerCell.Value2 = "right"
    } elseif ($answer -eq "b") {
        $answerCell.Value2 = "center"
    }
}

Write-Host "Updated face->book stimuli and expanded correct_ans codes through row $lastRow"
